$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (old column N), pushing the
# existing Late/heading/Outstanding columns one to the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# The freshly inserted column inherits the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab and move the selection.
$ws.Activate()
$ws.Range("S11").Select() | Out-Null
